$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Label" header in H1, matching the bold/bordered header
# style already used by the other header cells (B1:G1).
$ws.Range("H1").Value = "Label"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Fill in the new "Label" column values for both blocks of data rows
# (2-11 for the 100-iteration block, 12-21 for the 200-iteration block).
$labels = @(0, 0, 0, 0, 0, 1, 1, 1, 1, 1)

for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(2 + $i, 8).Value = $labels[$i]
    $ws.Cells.Item(12 + $i, 8).Value = $labels[$i]
}
